$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the report's column headers (row 1)
$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Curve"
$ws.Range("D1").Value = "Type"

# Correct the data point that no longer matches (row 2, column B)
$ws.Range("B2").Value = 637.4

# Leave the header row selected, as in the saved file
$ws.Range("A1:D1").Select()
